$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a "test_metadata" column just before the trailing "id" column on
#    the three sheets that carry the common Test*Data/Asset footer columns
#    (id, name, description, tags).
# ---------------------------------------------------------------------------

# TestAsset: id currently in column M (13)
$wsAsset = $wb.Worksheets.Item("TestAsset")
$wsAsset.Columns.Item(13).Insert()
$wsAsset.Cells.Item(1, 13).Value = "test_metadata"
$wsAsset.Range("F2:F1048576").Validation.Modify(3, 1, 1, '"Acceptable,BadButForgivable,NeverShow,TopAnswer"')

# AcceptanceTestAsset: id currently in column W (23)
$wsAccAsset = $wb.Worksheets.Item("AcceptanceTestAsset")
$wsAccAsset.Columns.Item(23).Insert()
$wsAccAsset.Cells.Item(1, 23).Value = "test_metadata"
$wsAccAsset.Range("P2:P1048576").Validation.Modify(3, 1, 1, '"Acceptable,BadButForgivable,NeverShow,TopAnswer"')

# TestEdgeData: id currently in column M (13)
$wsEdge = $wb.Worksheets.Item("TestEdgeData")
$wsEdge.Columns.Item(13).Insert()
$wsEdge.Cells.Item(1, 13).Value = "test_metadata"
$wsEdge.Range("F2:F1048576").Validation.Modify(3, 1, 1, '"Acceptable,BadButForgivable,NeverShow,TopAnswer"')

# ---------------------------------------------------------------------------
# 2. Append two new sheets: TestOutput and TestResultPKSet
#    (re-fetch the "after" target by name right before each Move so the
#    reference can't go stale across the intervening Add() calls)
# ---------------------------------------------------------------------------

$wsTestOutputNew = $wb.Worksheets.Add()
$wsTestOutputNew.Name = "TestOutput"
$afterTarget1 = $wb.Worksheets.Item("TestRunSession")
$wsTestOutputNew.Move($null, $afterTarget1)

# Object references go stale across Move() - re-resolve by name afterwards.
$wsTestOutput = $wb.Worksheets.Item("TestOutput")
$testOutputHeaders = @("test_suite_id", "test_case", "pks", "id", "name", "description", "tags")
for ($i = 0; $i -lt $testOutputHeaders.Length; $i++) {
    $wsTestOutput.Cells.Item(1, $i + 1).Value = $testOutputHeaders[$i]
}

$wsPkSetNew = $wb.Worksheets.Add()
$wsPkSetNew.Name = "TestResultPKSet"
$afterTarget2 = $wb.Worksheets.Item("TestOutput")
$wsPkSetNew.Move($null, $afterTarget2)

# Same re-resolve here.
$wsPkSet = $wb.Worksheets.Item("TestResultPKSet")
$pkSetHeaders = @("parent_pk", "merged_pk", "aragorn", "arax", "unsecret", "bte", "improving", "id", "name", "description", "tags")
for ($i = 0; $i -lt $pkSetHeaders.Length; $i++) {
    $wsPkSet.Cells.Item(1, $i + 1).Value = $pkSetHeaders[$i]
}
